# Adds maintenance log rows 7-12: six new equipment-1306 "Refrigeracion Liquida"
# entries, mirroring the existing row 6 record (workshop intake/repair log).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "1306"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = 45540
$ws.Range("B7").NumberFormat = "yyyy-mm-dd"
$ws.Range("C7").Value = 45540
$ws.Range("C7").NumberFormat = "yyyy-mm-dd"
$ws.Range("D7").Value = "cesar ramirez"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5587964476"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "UCL"
$ws.Range("G7").Value = "Refrigeracion Liquida"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "14123"
$ws.Range("H7").Style = "Normal"
$ws.Range("I7").Value = "NZXT"
$ws.Range("J7").NumberFormat = "@"
$ws.Range("J7").Value = "123412"
$ws.Range("J7").Style = "Normal"
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = "54324"
$ws.Range("K7").Style = "Normal"
$ws.Range("L7").Value = "No"
$ws.Range("M7").Value = "No"
$ws.Range("N7").Value = "No"
$ws.Range("O7").Value = "No"
$ws.Range("P7").Value = "Sí"
$ws.Range("Q7").Value = "Sí"
$ws.Range("R7").Value = "El equipo presenta falla por oxidacion y fuga de liquido"
$ws.Range("S7").Value = "Correctivo, Otro"
$ws.Range("T7").Value = "Se procedio a reparar la carcasa y sellar los tubos del radiador"
$ws.Range("U7").Value = "Si"
$ws.Range("V7").Value = "Alcohol Isopropílico, Aislantes, Liquido Limpiador Multiusos"
$ws.Range("W7").Value = "Juan Daniel Ramírez Zamora"
$ws.Range("X7").Value = "cesar ramirez"

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "1306"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = 45540
$ws.Range("B8").NumberFormat = "yyyy-mm-dd"
$ws.Range("C8").Value = 45540
$ws.Range("C8").NumberFormat = "yyyy-mm-dd"
$ws.Range("D8").Value = "cesar ramirez"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5587964476"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "UCL"
$ws.Range("G8").Value = "Refrigeracion Liquida"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "14123"
$ws.Range("H8").Style = "Normal"
$ws.Range("I8").Value = "NZXT"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "123412"
$ws.Range("J8").Style = "Normal"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "54324"
$ws.Range("K8").Style = "Normal"
$ws.Range("L8").Value = "No"
$ws.Range("M8").Value = "No"
$ws.Range("N8").Value = "No"
$ws.Range("O8").Value = "No"
$ws.Range("P8").Value = "Sí"
$ws.Range("Q8").Value = "Sí"
$ws.Range("R8").Value = "El equipo presenta falla por oxidacion y fuga de liquido"
$ws.Range("S8").Value = "Correctivo, Otro"
$ws.Range("T8").Value = "Se procedio a reparar la carcasa y sellar los tubos del radiador"
$ws.Range("U8").Value = "Si"
$ws.Range("V8").Value = "Alcohol Isopropílico, Aislantes, Liquido Limpiador Multiusos"
$ws.Range("W8").Value = "Juan Daniel Ramírez Zamora"

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "1306"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = 45540
$ws.Range("B9").NumberFormat = "yyyy-mm-dd"
$ws.Range("C9").Value = 45540
$ws.Range("C9").NumberFormat = "yyyy-mm-dd"
$ws.Range("D9").Value = "cesar ramirez"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5587964476"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "UCL"
$ws.Range("G9").Value = "Refrigeracion Liquida"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "14123"
$ws.Range("H9").Style = "Normal"
$ws.Range("I9").Value = "NZXT"
$ws.Range("J9").NumberFormat = "@"
$ws.Range("J9").Value = "123412"
$ws.Range("J9").Style = "Normal"
$ws.Range("K9").NumberFormat = "@"
$ws.Range("K9").Value = "54324"
$ws.Range("K9").Style = "Normal"
$ws.Range("L9").Value = "No"
$ws.Range("M9").Value = "No"
$ws.Range("N9").Value = "No"
$ws.Range("O9").Value = "No"
$ws.Range("P9").Value = "Sí"
$ws.Range("Q9").Value = "Sí"
$ws.Range("R9").Value = "El equipo presenta falla por oxidacion y fuga de liquido"
$ws.Range("S9").Value = "Correctivo, Otro"
$ws.Range("T9").Value = "Se procedio a reparar la carcasa y sellar los tubos del radiador"
$ws.Range("U9").Value = "Si"
$ws.Range("V9").Value = "Alcohol Isopropílico, Aislantes, Liquido Limpiador Multiusos"
$ws.Range("W9").Value = "Juan Daniel Ramírez Zamora"

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "1306"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = 45540
$ws.Range("B10").NumberFormat = "yyyy-mm-dd"
$ws.Range("C10").Value = 45540
$ws.Range("C10").NumberFormat = "yyyy-mm-dd"
$ws.Range("D10").Value = "cesar ramirez"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5587964476"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "UCL"
$ws.Range("G10").Value = "Refrigeracion Liquida"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "14123"
$ws.Range("H10").Style = "Normal"
$ws.Range("I10").Value = "NZXT"
$ws.Range("J10").NumberFormat = "@"
$ws.Range("J10").Value = "123412"
$ws.Range("J10").Style = "Normal"
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = "54324"
$ws.Range("K10").Style = "Normal"
$ws.Range("L10").Value = "No"
$ws.Range("M10").Value = "No"
$ws.Range("N10").Value = "No"
$ws.Range("O10").Value = "No"
$ws.Range("P10").Value = "Sí"
$ws.Range("Q10").Value = "Sí"
$ws.Range("R10").Value = "El equipo presenta falla por oxidacion y fuga de liquido"
$ws.Range("S10").Value = "Correctivo, Otro"
$ws.Range("T10").Value = "Se procedio a reparar la carcasa y sellar los tubos del radiador"
$ws.Range("U10").Value = "Si"
$ws.Range("V10").Value = "Alcohol Isopropílico, Aislantes, Liquido Limpiador Multiusos"
$ws.Range("W10").Value = "Juan Daniel Ramírez Zamora"

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "1306"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = 45540
$ws.Range("B11").NumberFormat = "yyyy-mm-dd"
$ws.Range("C11").Value = 45540
$ws.Range("C11").NumberFormat = "yyyy-mm-dd"
$ws.Range("D11").Value = "cesar ramirez"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5587964476"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "UCL"
$ws.Range("G11").Value = "Refrigeracion Liquida"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "14123"
$ws.Range("H11").Style = "Normal"
$ws.Range("I11").Value = "NZXT"
$ws.Range("J11").NumberFormat = "@"
$ws.Range("J11").Value = "123412"
$ws.Range("J11").Style = "Normal"
$ws.Range("K11").NumberFormat = "@"
$ws.Range("K11").Value = "54324"
$ws.Range("K11").Style = "Normal"
$ws.Range("L11").Value = "No"
$ws.Range("M11").Value = "No"
$ws.Range("N11").Value = "No"
$ws.Range("O11").Value = "No"
$ws.Range("P11").Value = "Sí"
$ws.Range("Q11").Value = "Sí"
$ws.Range("R11").Value = "El equipo presenta falla por oxidacion y fuga de liquido"
$ws.Range("S11").Value = "Correctivo, Otro"
$ws.Range("T11").Value = "Se procedio a reparar la carcasa y sellar los tubos del radiador"
$ws.Range("U11").Value = "No"
$ws.Range("V11").Value = "Alcohol Isopropílico, Aislantes, Liquido Limpiador Multiusos"
$ws.Range("W11").Value = "Juan Daniel Ramírez Zamora"
$ws.Range("X11").Value = "cesar ramirez"

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "1306"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = 45540
$ws.Range("B12").NumberFormat = "yyyy-mm-dd"
$ws.Range("C12").Value = 45540
$ws.Range("C12").NumberFormat = "yyyy-mm-dd"
$ws.Range("D12").Value = "cesar ramirez"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "5587964476"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "UCL"
$ws.Range("G12").Value = "Refrigeracion Liquida"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "14123"
$ws.Range("H12").Style = "Normal"
$ws.Range("I12").Value = "NZXT"
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "123412"
$ws.Range("J12").Style = "Normal"
$ws.Range("K12").NumberFormat = "@"
$ws.Range("K12").Value = "54324"
$ws.Range("K12").Style = "Normal"
$ws.Range("L12").Value = "No"
$ws.Range("M12").Value = "No"
$ws.Range("N12").Value = "No"
$ws.Range("O12").Value = "No"
$ws.Range("P12").Value = "Sí"
$ws.Range("Q12").Value = "Sí"
$ws.Range("R12").Value = "El equipo presenta falla por oxidacion y fuga de liquido"
$ws.Range("S12").Value = "Correctivo, Otro"
$ws.Range("T12").Value = "Se procedio a reparar la carcasa y sellar los tubos del radiador"
$ws.Range("U12").Value = "No"
$ws.Range("V12").Value = "Alcohol Isopropílico, Aislantes, Liquido Limpiador Multiusos"
$ws.Range("W12").Value = "Juan Daniel Ramírez Zamora"
$ws.Range("X12").Value = "cesar ramirez"
